$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) Create the new "2022-Q1" sheet by cloning the "2021-Q4" sheet
#    (same column layout / header / styles), placed right before the
#    "总计" (summary) sheet, then overwrite its data with the 2022-Q1
#    fund-holding figures and trim the two extra rows it doesn't need.
#
#    NOTE: worksheet object references in this host resolve
#    positionally, so any sheet handle fetched before a Copy/Add/
#    Delete that shifts tab order must be re-fetched afterwards
#    (by name) rather than reused.
# ------------------------------------------------------------------
$q4 = $wb.Worksheets.Item("2021-Q4")
$q4.Copy($wb.Worksheets.Item("总计"))

$q1 = $wb.Worksheets.Item("2021-Q4 (2)")
$q1.Name = "2022-Q1"

# Drop the two surplus data rows (source sheet has 8 rows, target has 6)
$q1.Rows.Item(8).Delete()
$q1.Rows.Item(8).Delete()

# Fund code (B) and numeric-looking figures (D:G) must stay text so
# leading zeros / trailing zeros survive - format the whole block as
# Text up front instead of per cell.
$q1.Range("B2:B7").NumberFormat = "@"
$q1.Range("D2:G7").NumberFormat = "@"

# Row 2
$q1.Cells.Item(2,1).Value = 0
$q1.Cells.Item(2,2).Value = "007012"
$q1.Cells.Item(2,3).Value = "湘财长顺混合A"
$q1.Cells.Item(2,4).Value = "4.70"
$q1.Cells.Item(2,5).Value = "94.08"
$q1.Cells.Item(2,6).Value = "6.07"
$q1.Cells.Item(2,7).Value = "0.2853"
$q1.Cells.Item(2,8).Value = 10

# Row 3
$q1.Cells.Item(3,1).Value = 1
$q1.Cells.Item(3,2).Value = "008128"
$q1.Cells.Item(3,3).Value = "湘财长源股票A"
$q1.Cells.Item(3,4).Value = "2.74"
$q1.Cells.Item(3,5).Value = "94.29"
$q1.Cells.Item(3,6).Value = "6.27"
$q1.Cells.Item(3,7).Value = "0.1718"
$q1.Cells.Item(3,8).Value = 7

# Row 4
$q1.Cells.Item(4,1).Value = 2
$q1.Cells.Item(4,2).Value = "007013"
$q1.Cells.Item(4,3).Value = "湘财长顺混合C"
$q1.Cells.Item(4,4).Value = "2.47"
$q1.Cells.Item(4,5).Value = "94.08"
$q1.Cells.Item(4,6).Value = "6.07"
$q1.Cells.Item(4,7).Value = "0.1499"
$q1.Cells.Item(4,8).Value = 10

# Row 5
$q1.Cells.Item(5,1).Value = 3
$q1.Cells.Item(5,2).Value = "011550"
$q1.Cells.Item(5,3).Value = "湘财创新成长一年持有期混合A"
$q1.Cells.Item(5,4).Value = "2.62"
$q1.Cells.Item(5,5).Value = "93.51"
$q1.Cells.Item(5,6).Value = "3.97"
$q1.Cells.Item(5,7).Value = "0.1040"
$q1.Cells.Item(5,8).Value = 9

# Row 6
$q1.Cells.Item(6,1).Value = 4
$q1.Cells.Item(6,2).Value = "008129"
$q1.Cells.Item(6,3).Value = "湘财长源股票C"
$q1.Cells.Item(6,4).Value = "1.05"
$q1.Cells.Item(6,5).Value = "94.29"
$q1.Cells.Item(6,6).Value = "6.27"
$q1.Cells.Item(6,7).Value = "0.0658"
$q1.Cells.Item(6,8).Value = 7

# Row 7
$q1.Cells.Item(7,1).Value = 5
$q1.Cells.Item(7,2).Value = "011551"
$q1.Cells.Item(7,3).Value = "湘财创新成长一年持有期混合C"
$q1.Cells.Item(7,4).Value = "0.28"
$q1.Cells.Item(7,5).Value = "93.51"
$q1.Cells.Item(7,6).Value = "3.97"
$q1.Cells.Item(7,7).Value = "0.0111"
$q1.Cells.Item(7,8).Value = 9

# ------------------------------------------------------------------
# 2) Update the "总计" (summary) sheet: insert a new row for the
#    2022-Q1 totals above the existing 2021-Q4 row, pushing the rest
#    down, and fill in its values (copy style/format from the row
#    being pushed down so the new row matches the sheet's existing
#    look). Re-fetch the sheet by name since the Copy() above shifted
#    tab positions.
# ------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()
$total.Range("A3:D3").Copy($total.Range("A2:D2"))

$total.Cells.Item(2,1).Value = 0
$total.Cells.Item(2,2).Value = "2022-Q1"
$total.Cells.Item(2,3).Value = 6
$total.Cells.Item(2,4).Value = 0.79

$total.Cells.Item(3,1).Value = 1
$total.Cells.Item(4,1).Value = 2
